$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spelling error on Timeline: remove the stray duplicate "kian " entry
# in A4 (shared string no longer referenced anywhere else).
$ws.Range("A4").ClearContents()

# Move the active selection to E5, matching the saved cursor position.
$ws.Range("E5").Select()
